# Update gh-pages to output generated at 456a3b4
# Refreshes "want to go" counts (and a couple of derived fields) across the
# four sheets of the 北京-漫展信息 workbook, flips one event to
# cancelled/not-sellable, and swaps a cover-image URL that was re-uploaded.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 264
$ws1.Range("F4").Value  = 9746
$ws1.Range("F5").Value  = 664
$ws1.Range("F6").Value  = 173
$ws1.Range("F7").Value  = 340
$ws1.Range("F9").Value  = 429
$ws1.Range("F11").Value = 211
$ws1.Range("F12").Value = 473
$ws1.Range("F13").Value = 12374
$ws1.Range("F23").Value = 166
$ws1.Range("F25").Value = 2108
$ws1.Range("F26").Value = 83
$ws1.Range("F28").Value = 61
$ws1.Range("F30").Value = 1046
$ws1.Range("F31").Value = 4220
$ws1.Range("F32").Value = 3709
$ws1.Range("F33").Value = 687
$ws1.Range("F36").Value = 45
$ws1.Range("F39").Value = 779
$ws1.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202404/Lv3lJQKQ1714287707594.jpeg"
$ws1.Range("F41").Value = 123
$ws1.Range("F42").Value = 451
$ws1.Range("F43").Value = 579

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Event cancelled: name gets a (取消) suffix and the price column becomes
# the literal "不可售" (not sellable) string instead of a numeric price.
$ws2.Range("C7").Value = "北京·DragonBand七龙珠限定Live纪念演出（取消）"
$ws2.Range("G7").Value = "不可售"

$ws2.Range("F14").Value = 41

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 56

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 264
$ws4.Range("F5").Value  = 9746
$ws4.Range("F6").Value  = 664
$ws4.Range("F8").Value  = 173
$ws4.Range("F9").Value  = 340
$ws4.Range("F11").Value = 429
$ws4.Range("F13").Value = 211
$ws4.Range("F14").Value = 473
$ws4.Range("F15").Value = 12374
$ws4.Range("F18").Value = 56
$ws4.Range("F23").Value = 166
$ws4.Range("F25").Value = 2108
$ws4.Range("F26").Value = 83
$ws4.Range("F27").Value = 61
$ws4.Range("F29").Value = 1046
$ws4.Range("F30").Value = 4220
$ws4.Range("F31").Value = 3709
$ws4.Range("F32").Value = 687
$ws4.Range("F35").Value = 45
$ws4.Range("F38").Value = 779
$ws4.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202404/Lv3lJQKQ1714287707594.jpeg"
$ws4.Range("F40").Value = 123
$ws4.Range("F41").Value = 451
$ws4.Range("F43").Value = 579
